# Update the marksheet's correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct-answer count 3 -> 5
$ws.Range("B11").Value = 5

# Total row: total marks 48 -> 80
$ws.Range("B12").Value = 80

# Corr/total marks display text 45/84 -> 80/140
$ws.Range("E12").Value = "80/140"
